$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.821.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.40%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.894.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.16%  "

$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7998"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.52%  "

$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3167"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.91%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.46"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07050"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08068"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("E12").Value = "  +3.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.883.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.83%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.345"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.98%  "

$ws.Range("E15").Value = "  +0.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.823.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.92%  "

$ws.Range("E18").Value = "  -0.96%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007718"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.221"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +19.04%  "

$ws.Range("E22").Value = "  -0.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.145.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.16%  "

$ws.Range("E25").Value = "  +5.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.346"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.97"
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.056"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.397"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.15%  "

$ws.Range("E31").Value = "  +1.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.444"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05708"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.048"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.54%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.263"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7381"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.93%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9991"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.630"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.39%  "

$ws.Range("E39").Value = "  -0.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.780"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4409"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.818"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8422"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.12%  "

$ws.Range("E45").Value = "  -0.18%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.033.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.27%  "

$ws.Range("E48").Value = "  -0.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.963"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.431"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.80%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.036.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.75%  "
